$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 621, shifting existing rows 621-668 down to 622-669
$ws.Rows.Item(621).Insert()

# Populate the newly inserted row 621 with the new record's data
$ws.Range("A621").Value = 3
$ws.Range("B621").Value = "Femacal de La Calera"
$ws.Range("C621").Value = "Coquimbo"
$ws.Range("D621").Value = 45265
$ws.Range("E621").Value = 5
$ws.Range("F621").Value = 100114013
$ws.Range("G621").Value = "Zanahoria"
$ws.Range("H621").Value = "Sin especificar"
$ws.Range("I621").Value = "Primera"
$ws.Range("J621").Value = 340
$ws.Range("K621").Value = 5500
$ws.Range("L621").Value = 6000
$ws.Range("M621").Value = 5765
$ws.Range("N621").Value = '$/saco 20 kilos'
$ws.Range("O621").Value = "Provincia de Quillota"
$ws.Range("P621").Value = 288
$ws.Range("Q621").Value = 20
$ws.Range("R621").Value = "Hortaliza"

Write-Host "done"
